$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-last data row (58) had its phone number stored as text;
# normalize it to a real number now that another row follows it.
$ws.Cells.Item(58, 1).Value = 76442781

# Append the new payment record as row 59.
# Phone number is text (leading apostrophe keeps "76442781" as a string
# instead of being auto-coerced to a number); same trick gives a true
# empty-string text cell for the two blank columns. Re-applying the
# "Normal" cell style afterwards drops the quote-prefix formatting that
# the apostrophe entry would otherwise leave behind, so no stray style
# is introduced.
$ws.Cells.Item(59, 1).Value = "'76442781"
$ws.Cells.Item(59, 1).Style = "Normal"
$ws.Cells.Item(59, 2).Value = "'"
$ws.Cells.Item(59, 2).Style = "Normal"
$ws.Cells.Item(59, 3).Value = "Cash"
$ws.Cells.Item(59, 4).Value = "2025-08-18T18:06:34"
$ws.Cells.Item(59, 5).Value = 120
$ws.Cells.Item(59, 6).Value = "'"
$ws.Cells.Item(59, 6).Style = "Normal"
$ws.Cells.Item(59, 7).Value = 115
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 9).Value = 100
$ws.Cells.Item(59, 10).Value = 5
